# Refresh the cryptocurrency symbol table (price + 1h volume-change columns)
# with the latest scraped values. Cells are stored as text (e.g. "307.59",
# "-4.01%"), so values are written with a leading apostrophe to keep Excel
# from auto-converting them to numeric/percentage types, and the style is
# reset back to "Normal" immediately after so no stray number-format / style
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'307.59"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'-4.01%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'39.98"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'-6.22%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.024"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'-3.93%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.07675"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'-5.91%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'4.237"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-2.37%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'1.620"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'-10.64%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.8906"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-6.36%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.1004"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-9.56%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.1730"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'-6.34%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.09009"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-3.23%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.04380"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'-5.00%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.1056"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'-0.32%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.001271"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'-1.39%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.005820"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'0.63%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'3.355"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'-0.37%"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'0.02%"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'-0.07%"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'7.057"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-5.26%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.1341"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'-3.76%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'0.3158"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'18.70%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.04221"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'0.62%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.001194"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'-4.66%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.004064"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'-5.69%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.0001222"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'-6.39%"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'-0.29%"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.02351"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'-9.13%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.05186"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'-5.58%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.007963"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'2.39%"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'-5.19%"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.006552"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-0.84%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.001988"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'-6.66%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.008093"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-4.57%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.3051"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-11.08%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00006572"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'-5.84%"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.00000000751"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'-0.25%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.003404"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-1.82%"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'41.18%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.00002103"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'-0.25%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.0002003"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'-0.25%"
$c.Style = "Normal"
